$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cell "Save" in H1, matching the style of the other
# header cells (B1:G1) by copying G1's formatting.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new "Save" column values for rows 2-5 (plain numeric cells,
# same default style as the other numeric columns).
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
